$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 766314
$ws.Range("E2").Value = 1429190161
$ws.Range("C10").Value = 345536
$ws.Range("E10").Value = 1817715752
$ws.Range("C13").Value = 187836
$ws.Range("E13").Value = 1165237431
$ws.Range("C81").Value = 88352
$ws.Range("E81").Value = 499616617
$ws.Range("C88").Value = 71267
$ws.Range("E88").Value = 110296657
$ws.Range("C91").Value = 18858
$ws.Range("E91").Value = 75166047
$ws.Range("C121").Value = 1306141
$ws.Range("E121").Value = 2274607696
$ws.Range("C129").Value = 633400
$ws.Range("E129").Value = 3427728680
$ws.Range("C132").Value = 585694
$ws.Range("E132").Value = 3463052588
$ws.Range("C136").Value = 26682
$ws.Range("E136").Value = 143703192
$ws.Range("C151").Value = 39929
$ws.Range("E151").Value = 60367248
$ws.Range("C154").Value = 18454
$ws.Range("E154").Value = 72765909
$ws.Range("C156").Value = 12401
$ws.Range("E156").Value = 40224583
$ws.Range("C160").Value = 4231
$ws.Range("E160").Value = 12512655
$ws.Range("C177").Value = 6938
$ws.Range("E177").Value = 30796837
$ws.Range("C207").Value = 154661
$ws.Range("E207").Value = 753651608
